$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.937933333333332
$ws.Range("H2").Value = 26.8138
$ws.Range("I2").Value = 0.2302024600837126
$ws.Range("J2").Value = 0.2302024600837126
$ws.Range("M2").Value = 31.61061466666667
$ws.Range("N2").Value = 94.831844
$ws.Range("O2").Value = 0.8860472269592234
$ws.Range("P2").Value = 0.8860472269592234
$ws.Range("Q2").Value = 282.5335665163555
$ws.Range("R2").Value = 2542.8020986472
$ws.Range("S2").Value = 0.2039702513963649
$ws.Range("T2").Value = 0.2039702513963649

# Row 3
$ws.Range("G3").Value = 8.937933333333332
$ws.Range("H3").Value = 26.8138
$ws.Range("I3").Value = 0.2302024600837126
$ws.Range("J3").Value = 0.2302024600837126
$ws.Range("O3").Value = 0.04688826274109129
$ws.Range("P3").Value = 0.04688826274109129
$ws.Range("Q3").Value = 14.95124379031111
$ws.Range("R3").Value = 134.5611941128
$ws.Range("S3").Value = 0.0107937934320507
$ws.Range("T3").Value = 0.0107937934320507

# Row 4
$ws.Range("G4").Value = 8.937933333333332
$ws.Range("H4").Value = 26.8138
$ws.Range("I4").Value = 0.2302024600837126
$ws.Range("J4").Value = 0.2302024600837126
$ws.Range("M4").Value = 2.392593
$ws.Range("N4").Value = 7.177778999999999
$ws.Range("O4").Value = 0.06706451029968528
$ws.Range("P4").Value = 0.06706451029968527
$ws.Range("Q4").Value = 21.38483672779999
$ws.Range("R4").Value = 192.4635305501999
$ws.Range("S4").Value = 0.01543841525529704
$ws.Range("T4").Value = 0.01543841525529703

# Row 5
$ws.Range("I5").Value = 0.5278886986241245
$ws.Range("J5").Value = 0.5278886986241244
$ws.Range("M5").Value = 31.61061466666667
$ws.Range("N5").Value = 94.831844
$ws.Range("O5").Value = 0.8860472269592234
$ws.Range("P5").Value = 0.8860472269592234
$ws.Range("Q5").Value = 647.8917588096787
$ws.Range("R5").Value = 5831.025829287108
$ws.Range("S5").Value = 0.4677343175590187
$ws.Range("T5").Value = 0.4677343175590186

# Row 6
$ws.Range("I6").Value = 0.5278886986241245
$ws.Range("J6").Value = 0.5278886986241244
$ws.Range("O6").Value = 0.04688826274109129
$ws.Range("P6").Value = 0.04688826274109129
$ws.Range("S6").Value = 0.0247517839991407
$ws.Range("T6").Value = 0.0247517839991407

# Row 7
$ws.Range("I7").Value = 0.5278886986241245
$ws.Range("J7").Value = 0.5278886986241244
$ws.Range("M7").Value = 2.392593
$ws.Range("N7").Value = 7.177778999999999
$ws.Range("O7").Value = 0.06706451029968528
$ws.Range("P7").Value = 0.06706451029968527
$ws.Range("Q7").Value = 49.038631587267
$ws.Range("R7").Value = 441.3476842854029
$ws.Range("S7").Value = 0.03540259706596505
$ws.Range("T7").Value = 0.03540259706596504

# Row 8
$ws.Range("G8").Value = 9.392449999999998
$ws.Range("H8").Value = 28.17735
$ws.Range("I8").Value = 0.241908841292163
$ws.Range("J8").Value = 0.2419088412921629
$ws.Range("M8").Value = 31.61061466666667
$ws.Range("N8").Value = 94.831844
$ws.Range("O8").Value = 0.8860472269592234
$ws.Range("P8").Value = 0.8860472269592234
$ws.Range("Q8").Value = 296.9011177259333
$ws.Range("R8").Value = 2672.1100595334
$ws.Range("S8").Value = 0.2143426580038399
$ws.Range("T8").Value = 0.2143426580038399

# Row 9
$ws.Range("G9").Value = 9.392449999999998
$ws.Range("H9").Value = 28.17735
$ws.Range("I9").Value = 0.241908841292163
$ws.Range("J9").Value = 0.2419088412921629
$ws.Range("O9").Value = 0.04688826274109129
$ws.Range("P9").Value = 0.04688826274109129
$ws.Range("Q9").Value = 15.71155260406667
$ws.Range("R9").Value = 141.4039734366
$ws.Range("S9").Value = 0.01134268530989989
$ws.Range("T9").Value = 0.01134268530989989

# Row 10
$ws.Range("G10").Value = 9.392449999999998
$ws.Range("H10").Value = 28.17735
$ws.Range("I10").Value = 0.241908841292163
$ws.Range("J10").Value = 0.2419088412921629
$ws.Range("M10").Value = 2.392593
$ws.Range("N10").Value = 7.177778999999999
$ws.Range("O10").Value = 0.06706451029968528
$ws.Range("P10").Value = 0.06706451029968527
$ws.Range("Q10").Value = 22.47231012284999
$ws.Range("R10").Value = 202.25079110565
$ws.Range("S10").Value = 0.0162234979784232
$ws.Range("T10").Value = 0.0162234979784232
